$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 75, shifting existing rows 75-143 down to 77-145.
$ws.Rows("75:76").Insert()

# New row 75 (Camote, 1a (cosecha))
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44658
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112045
$ws.Range("G75").Value = "Zapallo"
$ws.Range("H75").Value = "Camote"
$ws.Range("I75").Value = "1a (cosecha)"
$ws.Range("J75").Value = 200
$ws.Range("K75").Value = 300
$ws.Range("L75").Value = 350
$ws.Range("M75").Value = 325
$ws.Range("N75").Value = "$/kilo (volumen en unidades)"
$ws.Range("O75").Value = "Región de O'Higgins"
$ws.Range("P75").Value = 325
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"

# New row 76 (Paine, 1a (cosecha))
$ws.Range("A76").Value = 7
$ws.Range("B76").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C76").Value = "Ñuble"
$ws.Range("D76").Value = 44658
$ws.Range("E76").Value = 16
$ws.Range("F76").Value = 100112045
$ws.Range("G76").Value = "Zapallo"
$ws.Range("H76").Value = "Paine"
$ws.Range("I76").Value = "1a (cosecha)"
$ws.Range("J76").Value = 200
$ws.Range("K76").Value = 200
$ws.Range("L76").Value = 250
$ws.Range("M76").Value = 225
$ws.Range("N76").Value = "$/kilo (volumen en unidades)"
$ws.Range("O76").Value = "Región de O'Higgins"
$ws.Range("P76").Value = 225
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = "Hortaliza"

# Match the date display style used by other date cells in column D
$ws.Range("D75:D76").NumberFormat = $ws.Range("D77").NumberFormat
